$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices used in each data row: B,C,D,E,F,I,J,K,L,M,N (G stays 1, H is empty)
$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @{
    2 = @(1.02, 1.032440098722358, 1.041423314811248, 0.992614727750844, 1.050363779391306, 1.040737021158761, 1.037570308140926, 1.044202871759524, 0.9955398523335997, 1.053118263255288, 1.016526723653798)
    3 = @(1.02, 1.033166757806067, 1.041985927008585, 0.9936372048519299, 1.051062159088028, 1.040904785049461, 1.03794039522276, 1.04457676004986, 0.9963617723202687, 1.053629368519858, 1.016648822824291)
    4 = @(1.02, 1.03363747613114, 1.042350397637276, 0.9942998659930998, 1.051514831784722, 1.041012269272987, 1.03817967986911, 1.044818408201965, 0.9968940712668347, 1.05396018401667, 1.01672776305186)
    5 = @(1.02, 1.033835488912423, 1.042503720346199, 0.994578699834602, 1.051705318527317, 1.041057198505249, 1.038280229167391, 1.044919928003589, 0.9971179600053012, 1.05409928009314, 1.016760933245)
    6 = @(1.02, 1.033868743272297, 1.042529469666164, 0.994625531979634, 1.051737312763945, 1.041064727215768, 1.038297109102705, 1.044936969532303, 0.9971555583673455, 1.054122636129875, 1.016766501698587)
    7 = @(1.02, 1.033640121507621, 1.042352445955316, 0.994303590798249, 1.051517376361098, 1.041012870631104, 1.038181023596048, 1.044819764987769, 0.9968970624462089, 1.053962042544043, 1.016728206338139)
    8 = @(1.02, 1.032685567416685, 1.041613363474641, 0.9929600610674297, 1.05059963840554, 1.040793938796764, 1.03769541860671, 1.044329286705528, 0.9958175282591056, 1.053290972674081, 1.016568001053927)
    9 = @(1.02, 1.031007602398689, 1.040314332635844, 0.9906006454969559, 1.048988496135564, 1.040399996908484, 1.036838351508552, 1.043462900439626, 0.9939188001724441, 1.052109277886509, 1.016285212676942)
    10 = @(1.02, 1.029891822822464, 1.039450667774269, 0.989033133672735, 1.047918588648315, 1.040131940285502, 1.036266132504752, 1.042883982654257, 0.9926553831429383, 1.051322137088822, 1.01609638729451)
    11 = @(1.02, 1.029409382538104, 1.039077275297373, 0.988355674866747, 1.047456328417402, 1.040014592807539, 1.036018172734634, 1.042633008244475, 0.9921088820399291, 1.05098147464545, 1.01601455844554)
    12 = @(1.02, 1.029230289902469, 1.038938669913827, 0.9881042295826724, 1.047284779454552, 1.03997081377132, 1.035926042977867, 1.042539741689341, 0.9919059725120875, 1.050854965536438, 1.015984154083766)
    13 = @(1.02, 1.029268700998269, 1.038968397169366, 0.9881581567098651, 1.047321570235174, 1.039980213149594, 1.035945806298872, 1.042559749626301, 0.9919494934313052, 1.05088210086909, 1.015990676344304)
    14 = @(1.02, 1.02939457648693, 1.039065816301469, 0.9883348863814464, 1.047442144949912, 1.040010977908881, 1.036010557792065, 1.042625299680546, 0.9920921077337197, 1.050971016779363, 1.0160120454023)
    15 = @(1.02, 1.02947214679196, 1.039125851326187, 0.9884438009545853, 1.047516455628309, 1.04002990781634, 1.036050449858863, 1.042665681514018, 0.9921799884222134, 1.051025804583546, 1.016025210341356)
    16 = @(1.02, 1.029923855711147, 1.039475461000801, 0.9890781214508737, 1.047949288956979, 1.04013970142119, 1.036282584993085, 1.042900632820824, 0.9926916645766087, 1.051344749571573, 1.016101816650314)
    17 = @(1.02, 1.030207389544894, 1.03969491867548, 0.989476357848556, 1.048221067720051, 1.040208230761786, 1.036428148617472, 1.043047932502028, 0.9930127773699352, 1.051544863413643, 1.016149852364809)
    18 = @(1.02, 1.030372837400046, 1.039822980458778, 0.9897087662937556, 1.048379689623432, 1.040248079474224, 1.036513035459894, 1.043133820948522, 0.9932001317071769, 1.05166160316638, 1.016177864364793)
    19 = @(1.02, 1.030429262207182, 1.039866655625438, 0.9897880325774034, 1.048433792136807, 1.040261645914674, 1.03654197659223, 1.043163101731569, 0.9932640239640975, 1.051701411187595, 1.016187414629896)
    20 = @(1.02, 1.03017696206498, 1.039671367153873, 0.9894336180360679, 1.048191898282873, 1.040200890958146, 1.03641253286573, 1.043032131625732, 0.9929783193494215, 1.051523391348061, 1.016144699247911)
    21 = @(1.02, 1.029357506318067, 1.039037126316145, 0.9882828385668249, 1.047406634402472, 1.040001923718006, 1.035991490810651, 1.042605998012003, 0.9920501090198102, 1.050944832469757, 1.01600575300541)
    22 = @(1.02, 1.028842903015152, 1.038638871068714, 0.9875604150241495, 1.046913806155583, 1.039875720751158, 1.035726612555423, 1.042337820029939, 0.9914670000341481, 1.050581232971485, 1.015918337346331)
    23 = @(1.02, 1.02911564440862, 1.03884994400765, 0.9879432794643023, 1.047174977813133, 1.039942727703798, 1.035867043522361, 1.042480009543054, 0.991776070289318, 1.050773967839764, 1.015964683065939)
    24 = @(1.02, 1.030190710725538, 1.039682008900483, 0.9894529299347244, 1.048205078393618, 1.040204207880547, 1.036419589007435, 1.043039271451284, 0.9929938892766442, 1.051533093598785, 1.016147027739521)
    25 = @(1.02, 1.031440900583187, 1.040649756720545, 0.9912096547607049, 1.049404286746614, 1.040502801198036, 1.037060077934446, 1.0436871221942, 0.9944092447426414, 1.052414666563164, 1.016358375072634)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $colIndexes.Count; $i++) {
        $ws.Cells.Item($r, $colIndexes[$i]).Value = $rowVals[$i]
    }
}

Write-Output "Applied vm_pu updates for rows 2-25"